$d = $word.ActiveDocument

# 1. Remove the stray "_GoBack" bookmark that currently sits right after "cmd"
$d.Bookmarks.Item("_GoBack").Delete()

# 2. Wrap the "cd \Haie555.github.io" paragraph's text in a new "OLE_LINK2"
#    bookmark (mirrors the existing OLE_LINK1 bookmark already wrapping the
#    "git commit" line).
$cdPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("cd \Haie555.github.io")) {
        $cdPara = $p
        break
    }
}
$cdRange = $d.Range($cdPara.Range.Start, $cdPara.Range.End)
$d.Bookmarks.Add("OLE_LINK2", $cdRange) | Out-Null

# 3. Split "git push origin main" into "git push origin" and " main", with a
#    fresh (collapsed) "_GoBack" bookmark marking the split point - this is
#    what Word leaves behind after the cursor was left there on the last edit.
$pushPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("git push origin")) {
        $pushPara = $p
        break
    }
}
$splitRange = $d.Range($pushPara.Range.Start, $pushPara.Range.End)
$splitRange.Find.Execute("git push origin", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$splitPos = $splitRange.End
$d.Bookmarks.Add("_GoBack", $d.Range($splitPos, $splitPos)) | Out-Null
